# Update evaluation_metrics cross-validation results
# (check for abbreviation before lemmatization and showing results
#  of all folds of cross validation in table)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (kNN)
$ws.Range("B2").Value = 0.6329215282544041
$ws.Range("C2").Value = 0.6985840470620888
$ws.Range("D2").Value = 0.6329215282544041
$ws.Range("E2").Value = 0.6288039736425304
$ws.Range("G2").Value = 0.7338240528281872
$ws.Range("I2").Value = 0.6839428848500868
$ws.Range("J2").Value = 0.7016472203157174
$ws.Range("K2").Value = 0.7274869660793057
$ws.Range("L2").Value = 0.7016472203157174
$ws.Range("M2").Value = 0.7003135356634751
$ws.Range("N2").Value = 0.7832532601235416
$ws.Range("O2").Value = 0.796381171728196
$ws.Range("P2").Value = 0.7832532601235416
$ws.Range("Q2").Value = 0.783075680822123
$ws.Range("R2").Value = 0.6975062914664838
$ws.Range("S2").Value = 0.7648018239388501
$ws.Range("T2").Value = 0.6975062914664838
$ws.Range("U2").Value = 0.7013670086775078
$ws.Range("V2").Value = 0.6868222374742622
$ws.Range("W2").Value = 0.768815616616766
$ws.Range("X2").Value = 0.6868222374742622
$ws.Range("Y2").Value = 0.6951014546866139

# Row 3 (SVM)
$ws.Range("B3").Value = 0.7961564859299932
$ws.Range("C3").Value = 0.8085907744490285
$ws.Range("D3").Value = 0.7961564859299932
$ws.Range("E3").Value = 0.7963449218203591
$ws.Range("F3").Value = 0.8068634179821551
$ws.Range("G3").Value = 0.8175143763988478
$ws.Range("H3").Value = 0.8068634179821551
$ws.Range("I3").Value = 0.808185876519412
$ws.Range("J3").Value = 0.8262182566918325
$ws.Range("K3").Value = 0.8332262284052634
$ws.Range("L3").Value = 0.8262182566918325
$ws.Range("M3").Value = 0.8270849809206856
$ws.Range("N3").Value = 0.8240677190574239
$ws.Range("O3").Value = 0.8329728774185808
$ws.Range("P3").Value = 0.8240677190574239
$ws.Range("Q3").Value = 0.8233156931576824
$ws.Range("R3").Value = 0.8476549988560971
$ws.Range("S3").Value = 0.8562525997398776
$ws.Range("T3").Value = 0.8476549988560971
$ws.Range("U3").Value = 0.8463943668722983
$ws.Range("V3").Value = 0.8454815831617479
$ws.Range("W3").Value = 0.8533336711173224
$ws.Range("X3").Value = 0.8454815831617479
$ws.Range("Y3").Value = 0.8441787138429866

# Row 4 (LR)
$ws.Range("B4").Value = 0.8153283001601466
$ws.Range("C4").Value = 0.8245916392539296
$ws.Range("D4").Value = 0.8153283001601466
$ws.Range("E4").Value = 0.8148748851036117
$ws.Range("F4").Value = 0.8626172500571953
$ws.Range("G4").Value = 0.8688658601901917
$ws.Range("H4").Value = 0.8626172500571953
$ws.Range("I4").Value = 0.8628879067293121
$ws.Range("J4").Value = 0.8175017158544955
$ws.Range("K4").Value = 0.8251771835858351
$ws.Range("L4").Value = 0.8175017158544955
$ws.Range("M4").Value = 0.8172278659785916
$ws.Range("N4").Value = 0.8369023106840541
$ws.Range("O4").Value = 0.8412694402505408
$ws.Range("P4").Value = 0.8369023106840541
$ws.Range("Q4").Value = 0.8363400693236258
$ws.Range("R4").Value = 0.8498055364905056
$ws.Range("S4").Value = 0.8552800003368471
$ws.Range("T4").Value = 0.8498055364905056
$ws.Range("U4").Value = 0.8493243825931355
$ws.Range("V4").Value = 0.8476549988560971
$ws.Range("W4").Value = 0.8535742086068367
$ws.Range("X4").Value = 0.8476549988560971
$ws.Range("Y4").Value = 0.8470536906356572

# Row 5 (NB)
$ws.Range("B5").Value = 0.8090368336765043
$ws.Range("C5").Value = 0.8196520381624799
$ws.Range("D5").Value = 0.8090368336765043
$ws.Range("E5").Value = 0.8068330082309039
$ws.Range("F5").Value = 0.8412033859528713
$ws.Range("G5").Value = 0.8502483255424167
$ws.Range("H5").Value = 0.8412033859528713
$ws.Range("I5").Value = 0.8410518091077789
$ws.Range("N5").Value = 0.8218714253031344
$ws.Range("O5").Value = 0.8374887335214452
$ws.Range("P5").Value = 0.8218714253031344
$ws.Range("Q5").Value = 0.8196037335241346
$ws.Range("R5").Value = 0.8476549988560971
$ws.Range("S5").Value = 0.8589138787570496
$ws.Range("T5").Value = 0.8476549988560971
$ws.Range("U5").Value = 0.8464793942384269

# Row 6 (Ensemble)
$ws.Range("B6").Value = 0.8283230382063602
$ws.Range("C6").Value = 0.8361488568353399
$ws.Range("D6").Value = 0.8283230382063602
$ws.Range("E6").Value = 0.828597206776956
$ws.Range("G6").Value = 0.8570809370709608
$ws.Range("I6").Value = 0.851800851703387
$ws.Range("J6").Value = 0.8218485472431938
$ws.Range("K6").Value = 0.8303924661336406
$ws.Range("L6").Value = 0.8218485472431938
$ws.Range("M6").Value = 0.8218750059341182
$ws.Range("N6").Value = 0.8411805078929306
$ws.Range("O6").Value = 0.8498613990654673
$ws.Range("P6").Value = 0.8411805078929306
$ws.Range("Q6").Value = 0.8406354002454777
$ws.Range("R6").Value = 0.8518874399450928
$ws.Range("S6").Value = 0.8566227897235636
$ws.Range("T6").Value = 0.8518874399450928
$ws.Range("U6").Value = 0.8514615914589644
$ws.Range("V6").Value = 0.8583390528483186
$ws.Range("W6").Value = 0.8627061716931111
$ws.Range("X6").Value = 0.8583390528483186
$ws.Range("Y6").Value = 0.8582577235879147
